# "data from old laptop"
# The task's academy.html URL moved from the old EC2 IP address to localhost.
# Update the verify-script cell (E2) before the URL cell (D2) so that the
# shared-string table append order matches the canonical workbook (the
# "verify url equals ..." string ends up immediately before the new URL
# string).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tasks")

$ws.Range("E2").Value = "verify url equals ""http://localhost/academy.html""`n"
$ws.Range("D2").Value = "http://localhost/academy.html"

# The author's last selection before saving was cell C4.
$ws.Range("C4").Select()
